# Update the acquisition timestamp (column A) for data rows 2-11 on the
# "ランサーズ" sheet from 2025-10-23 06:27:52 to 2025-10-23 06:35:52.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-23 06:35:52"

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
